$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, shifting existing rows 101-191 down to 102-192
$ws.Rows(101).Insert()

# Populate the newly inserted row 101 with the new data record
$ws.Cells.Item(101, 1).Value = 5
$ws.Cells.Item(101, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(101, 3).Value = "Maule"
$ws.Cells.Item(101, 4).Value = 44741
$ws.Cells.Item(101, 5).Value = 7
$ws.Cells.Item(101, 6).Value = 100112017
$ws.Cells.Item(101, 7).Value = "Apio"
$ws.Cells.Item(101, 8).Value = "Americana (o)"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 600
$ws.Cells.Item(101, 11).Value = 7000
$ws.Cells.Item(101, 12).Value = 7000
$ws.Cells.Item(101, 13).Value = 7000
$ws.Cells.Item(101, 14).Value = "$/docena de matas"
$ws.Cells.Item(101, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(101, 16).Value = 1167
$ws.Cells.Item(101, 17).Value = 6
$ws.Cells.Item(101, 18).Value = "Hortaliza"
